$d = $word.ActiveDocument

# Step 1: find the target run's text "+2 card on draw" (exact, case sensitive)
# and turn it into "+2 card on draw " (trailing space added), leaving its
# formatting (bold off, sz 16) untouched.
$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute("+2 card on draw", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "+2 card on draw ", 2)

# Step 2: collapse to the end of the replaced text and insert a new run
# "IMPLEMENTED" right after it, matching the formatting used elsewhere in
# the document for these markers (bold, sz16, lang en-US).
$rng.Collapse(0)
$rng.InsertAfter("IMPLEMENTED")
$rng.Font.Bold = $true
$rng.Font.Size = 16
$rng.LanguageID = "wdEnglishUS"
